# Prepare public release: split the single-sheet sample workbook into
# three sheets (Sheet1/Sheet2/Sheet3), each carrying its own small data
# sample, and shrink Sheet1's data down to a single header-ish row.

$wb = $excel.ActiveWorkbook

# --- Sheet1: rename the existing sheet and rewrite its contents ----------
$ws1 = $wb.ActiveSheet
$ws1.Name = "Sheet1"
$ws1.Range("A1").Value = "Data1"
$ws1.Range("B1").Value = "Value1"
# Row 2 ("Test"/"123") goes away entirely so the used range shrinks to A1:B1.
$ws1.Range("A2:B2").ClearContents()

# --- Sheet2: new sheet right after Sheet1 with its own sample row --------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("A1").Value = "Data2"
$ws2.Range("B1").Value = "Value2"

# --- Sheet3: new, empty sheet right after Sheet2 --------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"

# Leave the workbook focused back on the first sheet (matches activeTab=0).
$ws1.Activate()
